# Applies the OOXML diff:
#  - Slide 2 ("Nós operacionais"): clear the "-Produtos novos e atualizado "
#    line and reset the body placeholder's autofit.
#  - Slide 3 ("Nós operacionais" variant): remove the "-Manutenção de
#    defeitos" paragraph, reset autofit, and reflow the diagonal-corner
#    rectangle + its curved connector to the new (shorter) size.
#  - Slide 4 ("Nós operacionais" variant): reword two bullet lines.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - "Espaço Reservado para Texto 4" placeholder (shape 2)
# ---------------------------------------------------------------------
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tf2 = $sh2.TextFrame

# Reset the auto-shrunk text back to its natural (un-scaled) normAutofit.
$tf2.AutoSize = 2

# Remove the "-Produtos novos e atualizado " bullet text (3rd paragraph),
# leaving the now-empty paragraph in place.
$tf2.TextRange.Paragraphs(3).Text = ""

# ---------------------------------------------------------------------
# Slide 3 - "Espaço Reservado para Texto 4" placeholder (shape 2)
# ---------------------------------------------------------------------
$s3  = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tf3 = $sh3.TextFrame

# Reset the auto-shrunk text back to its natural (un-scaled) normAutofit.
$tf3.AutoSize = 2

# Drop the whole "-Manutenção de defeitos" paragraph (2nd paragraph).
$tf3.TextRange.Paragraphs(2).Delete()

# ---------------------------------------------------------------------
# Slide 3 - "Retângulo: Cantos Diagonais Arredondados 7" (shape 3) and its
# "Conector: Curvo 16" curved connector (shape 4) move up & shrink to
# match the now-shorter text placeholder.
# ---------------------------------------------------------------------
$rect = $s3.Shapes.Item(3)
$conn = $s3.Shapes.Item(4)

$emuPerPt = 12700
# Tiny epsilon nudges the point value onto the correct side of the
# engine's EMU rounding boundary (Shape.Top/Height are single-precision
# points, same as real PowerPoint).
$eps = 0.00001

$rect.Top    = (643624 / $emuPerPt) + $eps
$rect.Height = (1009940 / $emuPerPt) + $eps

$conn.Top    = (1148594 / $emuPerPt) + $eps
$conn.Height = (2146692 / $emuPerPt) + $eps
$conn.Adjustments(2) = 0.61762

# ---------------------------------------------------------------------
# Slide 4 - "Espaço Reservado para Texto 4" placeholder (shape 2): reword
# the two bullet lines (formatting/runs stay untouched).
# ---------------------------------------------------------------------
$s4  = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tf4 = $sh4.TextFrame
$tr4 = $tf4.TextRange

$tr4.Paragraphs(2).Runs(1).Text = "-Prestar manutenção "
$tr4.Paragraphs(3).Runs(1).Text = "-Reparo da rede"
